$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.040.88"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.565.37"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.45"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.95"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0596"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0860"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "1.788.70"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "1.566.61"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.77"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "27.004.25"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.96"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "0.0₃0704"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "215.18"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.17"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.82"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.61"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0473"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.13"
$ws.Range("E31").Value = "  +4.43%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.19"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("D34").Value = "1.430.85"
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("E35").Value = "  +19.57%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.530"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.83"
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.38"
$ws.Range("E42").Value = "  +4.35%  "
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.49"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "1.705.44"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.23"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0517"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0962"
$ws.Range("E51").Value = "  +0.56%  "
